# Generate Report for Handback
# Refresh the timestamp columns ("Latest HO Xliff Generate Date",
# "Correspond Handoff Datetime", "Correspond Handback DateTime") on each
# sheet to reflect the new report-generation run, preserving their
# existing date/time display format.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$dateFormat = "yyyy-mm-dd HH:mm:ss"

# Overview sheet: Latest HO Xliff Generate Date
$overview.Range("G2").NumberFormat = $dateFormat
$overview.Range("G2").Value = "2016-08-16 18:59:21"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$zhcn.Range("H2").NumberFormat = $dateFormat
$zhcn.Range("H2").Value = "2016-08-16 18:59:16"
$zhcn.Range("K2").NumberFormat = $dateFormat
$zhcn.Range("K2").Value = "2016-08-16 18:59:36"

# de-de sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$dede.Range("H2").NumberFormat = $dateFormat
$dede.Range("H2").Value = "2016-08-16 18:59:21"
$dede.Range("K2").NumberFormat = $dateFormat
$dede.Range("K2").Value = "2016-08-16 18:59:44"
